$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 16.23365402659878
$ws.Cells.Item(2, 3).Value = 3.86408143130137
$ws.Cells.Item(2, 4).Value = 13.45272583510759
$ws.Cells.Item(2, 5).Value = 13.43235276683116
$ws.Cells.Item(2, 7).Value = 68.88191477822151
$ws.Cells.Item(2, 8).Value = 24.30064963152257
$ws.Cells.Item(2, 9).Value = 39.08158368773766
$ws.Cells.Item(2, 10).Value = 8.211349561222056
$ws.Cells.Item(2, 11).Value = 14.01310740033882
$ws.Cells.Item(2, 12).Value = 13.10723826090272
$ws.Cells.Item(2, 13).Value = 17.97745759718419

$ws.Cells.Item(3, 2).Value = 16.16842270184804
$ws.Cells.Item(3, 3).Value = 3.681421306795692
$ws.Cells.Item(3, 4).Value = 13.45544950731072
$ws.Cells.Item(3, 5).Value = 13.4523381950136
$ws.Cells.Item(3, 7).Value = 68.75018976283606
$ws.Cells.Item(3, 8).Value = 24.31205263255582
$ws.Cells.Item(3, 9).Value = 39.08417458963331
$ws.Cells.Item(3, 10).Value = 8.211500009759707
$ws.Cells.Item(3, 11).Value = 13.97276531667602
$ws.Cells.Item(3, 12).Value = 13.12495509120098
$ws.Cells.Item(3, 13).Value = 17.99322053149979

$ws.Cells.Item(4, 2).Value = 16.13244529166927
$ws.Cells.Item(4, 3).Value = 3.563592435552464
$ws.Cells.Item(4, 4).Value = 13.45924722181303
$ws.Cells.Item(4, 5).Value = 13.46580756648313
$ws.Cells.Item(4, 7).Value = 68.67998569267961
$ws.Cells.Item(4, 8).Value = 24.32191551263523
$ws.Cells.Item(4, 9).Value = 39.09049452836956
$ws.Cells.Item(4, 10).Value = 8.211620046196519
$ws.Cells.Item(4, 11).Value = 13.95137258553986
$ws.Cells.Item(4, 12).Value = 13.13745921237892
$ws.Cells.Item(4, 13).Value = 18.00580830433204

$ws.Cells.Item(5, 2).Value = 16.11882232567873
$ws.Cells.Item(5, 3).Value = 3.514175684209019
$ws.Cells.Item(5, 4).Value = 13.46133003658027
$ws.Cells.Item(5, 5).Value = 13.47159822616512
$ws.Cells.Item(5, 7).Value = 68.65407351516295
$ws.Cells.Item(5, 8).Value = 24.32665445977838
$ws.Cells.Item(5, 9).Value = 39.09425972148824
$ws.Cells.Item(5, 10).Value = 8.211675920557948
$ws.Cells.Item(5, 11).Value = 13.94351186609724
$ws.Cells.Item(5, 12).Value = 13.14296396531469
$ws.Cells.Item(5, 13).Value = 18.01167034728919

$ws.Cells.Item(6, 2).Value = 16.11662330911243
$ws.Cells.Item(6, 3).Value = 3.505886350898844
$ws.Cells.Item(6, 4).Value = 13.46170823119595
$ws.Cells.Item(6, 5).Value = 13.47257800147447
$ws.Cells.Item(6, 7).Value = 68.64993399661752
$ws.Cells.Item(6, 8).Value = 24.32748483276131
$ws.Cells.Item(6, 9).Value = 39.09495679879149
$ws.Cells.Item(6, 10).Value = 8.211685618827866
$ws.Cells.Item(6, 11).Value = 13.94225856667034
$ws.Cells.Item(6, 12).Value = 13.14390275151849
$ws.Cells.Item(6, 13).Value = 18.01268799352302

$ws.Cells.Item(7, 2).Value = 16.13225734751692
$ws.Cells.Item(7, 3).Value = 3.562931616130691
$ws.Cells.Item(7, 4).Value = 13.45927314342014
$ws.Cells.Item(7, 5).Value = 13.46588443881449
$ws.Cells.Item(7, 7).Value = 68.67962529941019
$ws.Cells.Item(7, 8).Value = 24.32197650938063
$ws.Cells.Item(7, 9).Value = 39.09054048928798
$ws.Cells.Item(7, 10).Value = 8.211620771557232
$ws.Cells.Item(7, 11).Value = 13.95126309372892
$ws.Cells.Item(7, 12).Value = 13.13753179401426
$ws.Cells.Item(7, 13).Value = 18.00588439533269

$ws.Cells.Item(8, 2).Value = 16.21032500401286
$ws.Cells.Item(8, 3).Value = 3.802292436424528
$ws.Cells.Item(8, 4).Value = 13.45322435070225
$ws.Cells.Item(8, 5).Value = 13.43899533469819
$ws.Cells.Item(8, 7).Value = 68.83428808107142
$ws.Cells.Item(8, 8).Value = 24.30398767481617
$ws.Cells.Item(8, 9).Value = 39.0814959431008
$ws.Cells.Item(8, 10).Value = 8.211395696657545
$ws.Cells.Item(8, 11).Value = 13.99850135131538
$ws.Cells.Item(8, 12).Value = 13.11300979818566
$ws.Cells.Item(8, 13).Value = 17.98228940429881

$ws.Cells.Item(9, 2).Value = 16.39509444835664
$ws.Cells.Item(9, 3).Value = 4.225718431087071
$ws.Cells.Item(9, 4).Value = 13.4581867729159
$ws.Cells.Item(9, 5).Value = 13.39575247078823
$ws.Cells.Item(9, 7).Value = 69.22168137618922
$ws.Cells.Item(9, 8).Value = 24.29140286956577
$ws.Cells.Item(9, 9).Value = 39.10124918957667
$ws.Cells.Item(9, 10).Value = 8.211173721480245
$ws.Cells.Item(9, 11).Value = 14.11755439402634
$ws.Cells.Item(9, 12).Value = 13.07780621976776
$ws.Cells.Item(9, 13).Value = 17.95906028037157

$ws.Cells.Item(10, 2).Value = 16.54919618694752
$ws.Cells.Item(10, 3).Value = 4.507814448541621
$ws.Cells.Item(10, 4).Value = 13.47202715549956
$ws.Cells.Item(10, 5).Value = 13.36973707172501
$ws.Cells.Item(10, 7).Value = 69.55664757384424
$ws.Cells.Item(10, 8).Value = 24.29596882994611
$ws.Cells.Item(10, 9).Value = 39.1385595113406
$ws.Cells.Item(10, 10).Value = 8.211144370695925
$ws.Cells.Item(10, 11).Value = 14.22056122263894
$ws.Cells.Item(10, 12).Value = 13.05977291102024
$ws.Cells.Item(10, 13).Value = 17.95597426159139

$ws.Cells.Item(11, 2).Value = 16.62305689729705
$ws.Cells.Item(11, 3).Value = 4.629708640932777
$ws.Cells.Item(11, 4).Value = 13.48052014758939
$ws.Cells.Item(11, 5).Value = 13.35914560206222
$ws.Cells.Item(11, 7).Value = 69.71972844186415
$ws.Cells.Item(11, 8).Value = 24.3010376993628
$ws.Cells.Item(11, 9).Value = 39.16046431995849
$ws.Cells.Item(11, 10).Value = 8.211160064218641
$ws.Cells.Item(11, 11).Value = 14.27065718890854
$ws.Cells.Item(11, 12).Value = 13.05326367136116
$ws.Cells.Item(11, 13).Value = 17.95758837132841

$ws.Cells.Item(12, 2).Value = 16.65154570916211
$ws.Cells.Item(12, 3).Value = 4.674933009882831
$ws.Cells.Item(12, 4).Value = 13.48405033074036
$ws.Cells.Item(12, 5).Value = 13.35531313530376
$ws.Cells.Item(12, 7).Value = 69.78299850939422
$ws.Cells.Item(12, 8).Value = 24.30338637459293
$ws.Cells.Item(12, 9).Value = 39.1694659594915
$ws.Cells.Item(12, 10).Value = 8.211170183289491
$ws.Cells.Item(12, 11).Value = 14.29007987938989
$ws.Cells.Item(12, 12).Value = 13.05104184091009
$ws.Cells.Item(12, 13).Value = 17.95863170856973

$ws.Cells.Item(13, 2).Value = 16.64538738183257
$ws.Cells.Item(13, 3).Value = 4.665234812918808
$ws.Cells.Item(13, 4).Value = 13.48327610572866
$ws.Cells.Item(13, 5).Value = 13.35613060334423
$ws.Cells.Item(13, 7).Value = 69.7693051989674
$ws.Cells.Item(13, 8).Value = 24.30286147428503
$ws.Cells.Item(13, 9).Value = 39.16749590520918
$ws.Cells.Item(13, 10).Value = 8.21116781823017
$ws.Cells.Item(13, 11).Value = 14.28587695010695
$ws.Cells.Item(13, 12).Value = 13.05150955013798
$ws.Cells.Item(13, 13).Value = 17.95838782097343

$ws.Cells.Item(14, 2).Value = 16.62539040472402
$ws.Cells.Item(14, 3).Value = 4.633448047955331
$ws.Cells.Item(14, 4).Value = 13.48080429469956
$ws.Cells.Item(14, 5).Value = 13.35882673211055
$ws.Cells.Item(14, 7).Value = 69.72490349326681
$ws.Cells.Item(14, 8).Value = 24.30122233079416
$ws.Cells.Item(14, 9).Value = 39.16119074116553
$ws.Cells.Item(14, 10).Value = 8.211160813023545
$ws.Cells.Item(14, 11).Value = 14.27224611567169
$ws.Cells.Item(14, 12).Value = 13.05307601170289
$ws.Cells.Item(14, 13).Value = 17.95766555862786

$ws.Cells.Item(15, 2).Value = 16.61320867111542
$ws.Cells.Item(15, 3).Value = 4.613855820706897
$ws.Cells.Item(15, 4).Value = 13.47933108227649
$ws.Cells.Item(15, 5).Value = 13.36050139410466
$ws.Cells.Item(15, 7).Value = 69.69790269814334
$ws.Cells.Item(15, 8).Value = 24.30027417022474
$ws.Cells.Item(15, 9).Value = 39.15742060782289
$ws.Cells.Item(15, 10).Value = 8.211157065998359
$ws.Cells.Item(15, 11).Value = 14.2639553542194
$ws.Cells.Item(15, 12).Value = 13.05406715344246
$ws.Cells.Item(15, 13).Value = 17.95727936650761

$ws.Cells.Item(16, 2).Value = 16.5444430679217
$ws.Cells.Item(16, 3).Value = 4.499718090155847
$ws.Cells.Item(16, 4).Value = 13.47151617325358
$ws.Cells.Item(16, 5).Value = 13.37045422381759
$ws.Cells.Item(16, 7).Value = 69.54620320288024
$ws.Cells.Item(16, 8).Value = 24.29569768730893
$ws.Cells.Item(16, 9).Value = 39.13722702522875
$ws.Cells.Item(16, 10).Value = 8.211143929513618
$ws.Cells.Item(16, 11).Value = 14.2173513771695
$ws.Cells.Item(16, 12).Value = 13.06023235116051
$ws.Cells.Item(16, 13).Value = 17.95592934949335

$ws.Cells.Item(17, 2).Value = 16.50320560781133
$ws.Cells.Item(17, 3).Value = 4.428043413628933
$ws.Cells.Item(17, 4).Value = 13.46728341768455
$ws.Cells.Item(17, 5).Value = 13.37687801756688
$ws.Cells.Item(17, 7).Value = 69.45586517204944
$ws.Cells.Item(17, 8).Value = 24.29365576709078
$ws.Cells.Item(17, 9).Value = 39.1261005185666
$ws.Cells.Item(17, 10).Value = 8.211143309548467
$ws.Cells.Item(17, 11).Value = 14.18958147592511
$ws.Cells.Item(17, 12).Value = 13.0644480521371
$ws.Cells.Item(17, 13).Value = 17.95587286230284

$ws.Cells.Item(18, 2).Value = 16.47984206175246
$ws.Cells.Item(18, 3).Value = 4.386212967169356
$ws.Cells.Item(18, 4).Value = 13.46505575742739
$ws.Cells.Item(18, 5).Value = 13.38068984493889
$ws.Cells.Item(18, 7).Value = 69.40491412012616
$ws.Cells.Item(18, 8).Value = 24.29276298123677
$ws.Cells.Item(18, 9).Value = 39.12016525566217
$ws.Cells.Item(18, 10).Value = 8.211145686860471
$ws.Cells.Item(18, 11).Value = 14.17391454009538
$ws.Cells.Item(18, 12).Value = 13.06703231392606
$ws.Cells.Item(18, 13).Value = 17.95612456112158

$ws.Cells.Item(19, 2).Value = 16.47199319411045
$ws.Cells.Item(19, 3).Value = 4.371946294327856
$ws.Cells.Item(19, 4).Value = 13.46433709965999
$ws.Cells.Item(19, 5).Value = 13.38200057824571
$ws.Cells.Item(19, 7).Value = 69.38783697153137
$ws.Cells.Item(19, 8).Value = 24.29250910492444
$ws.Cells.Item(19, 9).Value = 39.11823550971463
$ws.Cells.Item(19, 10).Value = 8.211146961342424
$ws.Cells.Item(19, 11).Value = 14.16866285602472
$ws.Cells.Item(19, 12).Value = 13.06793471109586
$ws.Cells.Item(19, 13).Value = 17.95625864538564

$ws.Cells.Item(20, 2).Value = 16.50755879962926
$ws.Cells.Item(20, 3).Value = 4.435736001687425
$ws.Cells.Item(20, 4).Value = 13.46771260019279
$ws.Cells.Item(20, 5).Value = 13.37618208477539
$ws.Cells.Item(20, 7).Value = 69.46537755176165
$ws.Cells.Item(20, 8).Value = 24.2938439876642
$ws.Cells.Item(20, 9).Value = 39.12723690623469
$ws.Cells.Item(20, 10).Value = 8.211143092601299
$ws.Cells.Item(20, 11).Value = 14.19250609605799
$ws.Cells.Item(20, 12).Value = 13.06398277951354
$ws.Cells.Item(20, 13).Value = 17.95584947588586

$ws.Cells.Item(21, 2).Value = 16.63125007955758
$ws.Cells.Item(21, 3).Value = 4.642810012233446
$ws.Cells.Item(21, 4).Value = 13.48152181740249
$ws.Cells.Item(21, 5).Value = 13.35802997860738
$ws.Cells.Item(21, 7).Value = 69.73790443933862
$ws.Cells.Item(21, 8).Value = 24.30169214739016
$ws.Cells.Item(21, 9).Value = 39.16302356088676
$ws.Cells.Item(21, 10).Value = 8.211162757281745
$ws.Cells.Item(21, 11).Value = 14.27623765145953
$ws.Cells.Item(21, 12).Value = 13.05260931151565
$ws.Cells.Item(21, 13).Value = 17.95786599317013

$ws.Cells.Item(22, 2).Value = 16.71510629001835
$ws.Cells.Item(22, 3).Value = 4.772696029273754
$ws.Cells.Item(22, 4).Value = 13.49237659095373
$ws.Cells.Item(22, 5).Value = 13.34720553770462
$ws.Cells.Item(22, 7).Value = 69.92483386109231
$ws.Cells.Item(22, 8).Value = 24.30932251543443
$ws.Cells.Item(22, 9).Value = 39.19053035717338
$ws.Cells.Item(22, 10).Value = 8.211199951781509
$ws.Cells.Item(22, 11).Value = 14.33359161130179
$ws.Cells.Item(22, 12).Value = 13.04659264596012
$ws.Cells.Item(22, 13).Value = 17.96170155584488

$ws.Cells.Item(23, 2).Value = 16.6700818409981
$ws.Cells.Item(23, 3).Value = 4.703874533631261
$ws.Cells.Item(23, 4).Value = 13.48641644174722
$ws.Cells.Item(23, 5).Value = 13.35288782939993
$ws.Cells.Item(23, 7).Value = 69.82426780418716
$ws.Cells.Item(23, 8).Value = 24.30502156805316
$ws.Cells.Item(23, 9).Value = 39.17547356571422
$ws.Cells.Item(23, 10).Value = 8.211177873111387
$ws.Cells.Item(23, 11).Value = 14.30274455155099
$ws.Cells.Item(23, 12).Value = 13.04967443191825
$ws.Cells.Item(23, 13).Value = 17.9594247481451

$ws.Cells.Item(24, 2).Value = 16.50558964715978
$ws.Cells.Item(24, 3).Value = 4.43226012471407
$ws.Cells.Item(24, 4).Value = 13.46751792558546
$ws.Cells.Item(24, 5).Value = 13.37649634621244
$ws.Cells.Item(24, 7).Value = 69.46107393320197
$ws.Cells.Item(24, 8).Value = 24.29375801725256
$ws.Cells.Item(24, 9).Value = 39.12672170750069
$ws.Cells.Item(24, 10).Value = 8.211143182170197
$ws.Cells.Item(24, 11).Value = 14.1911829449564
$ws.Cells.Item(24, 12).Value = 13.06419262897098
$ws.Cells.Item(24, 13).Value = 17.9558591635562

$ws.Cells.Item(25, 2).Value = 16.34181234975138
$ws.Cells.Item(25, 3).Value = 4.116204687267362
$ws.Cells.Item(25, 4).Value = 13.45504879840865
$ws.Cells.Item(25, 5).Value = 13.4064380908526
$ws.Cells.Item(25, 7).Value = 69.10796485031754
$ws.Cells.Item(25, 8).Value = 24.29237969318426
$ws.Cells.Item(25, 9).Value = 39.09189636797955
$ws.Cells.Item(25, 10).Value = 8.211210286195485
$ws.Cells.Item(25, 11).Value = 14.08257414996768
$ws.Cells.Item(25, 12).Value = 13.08595267498205
$ws.Cells.Item(25, 13).Value = 17.9628837554385
